# Insert two new weekly records for "Femacal de La Calera - Alcachofa" at the
# top of the existing block of "Argentina(o)/Primera" and "Espanola/Primera"
# rows (old rows 555-556), pushing all subsequent rows down by two.
#
# The new rows reuse the same Variedad/Calidad/Unidad de comercializacion/
# Origen as the rows that used to occupy 555-556, only the date and the
# volume/price figures differ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 555.
$ws.Rows("555:556").Insert()

# The old row 555 / 556 are now at 557 / 558 - clone them into the two new
# blank rows so every column (Mercado, Region, Variedad, Calidad, Unidad,
# Origen, Clasificacion, ...) keeps its value, then overwrite the few cells
# that actually carry new data for this week.
$ws.Range("A557:R557").Copy()
$ws.Range("A555").PasteSpecial()

$ws.Range("A558:R558").Copy()
$ws.Range("A556").PasteSpecial()

$ws.Application.CutCopyMode = 0

# Row 555: Alcachofa / Argentina(o) / Primera
$ws.Cells.Item(555, 4).Value = 45166   # Fecha
$ws.Cells.Item(555, 10).Value = 50     # Volumen
$ws.Cells.Item(555, 11).Value = 13000  # Precio minimo
$ws.Cells.Item(555, 12).Value = 13000  # Precio maximo
$ws.Cells.Item(555, 13).Value = 13000  # Precio promedio ponderado
$ws.Cells.Item(555, 16).Value = 260    # Precio $/Kg

# Row 556: Alcachofa / Espanola / Primera
$ws.Cells.Item(556, 4).Value = 45166   # Fecha
$ws.Cells.Item(556, 10).Value = 55     # Volumen
$ws.Cells.Item(556, 11).Value = 14000  # Precio minimo
$ws.Cells.Item(556, 12).Value = 14000  # Precio maximo
$ws.Cells.Item(556, 13).Value = 14000  # Precio promedio ponderado
$ws.Cells.Item(556, 16).Value = 467    # Precio $/Kg
